$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.619.17'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '1.893.99'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.56'
$ws.Range("E5").Value = '  +1.29%  '
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("E7").Value = '  +0.56%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2941'
$ws.Range("E8").Value = '  +2.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06715'
$ws.Range("E9").Value = '  +0.83%  '
$ws.Range("D10").Value = '1.880.61'
$ws.Range("E10").Value = '  -0.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '17.16'
$ws.Range("E11").Value = '  +2.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07349'
$ws.Range("E12").Value = '  +1.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.166'
$ws.Range("E13").Value = '  +3.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.18'
$ws.Range("E14").Value = '  -0.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6694'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = '30.558.60'
$ws.Range("E16").Value = '  -0.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007872'
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("E18").Value = '  +3.75%  '
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").Value = '2.142.32'
$ws.Range("E20").Value = '  +0.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.298'
$ws.Range("E21").Value = '  +11.95%  '
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '190.60'
$ws.Range("E23").Value = '  +2.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.182'
$ws.Range("E24").Value = '  +2.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.478'
$ws.Range("E25").Value = '  +2.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.13'
$ws.Range("E26").Value = '  +2.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.34'
$ws.Range("E27").Value = '  +0.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.937'
$ws.Range("E28").Value = '  +5.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.477'
$ws.Range("E29").Value = '  +5.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.381'
$ws.Range("E30").Value = '  +2.96%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09164'
$ws.Range("E31").Value = '  +1.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.127'
$ws.Range("E32").Value = '  +4.99%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05226'
$ws.Range("E33").Value = '  +0.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7398'
$ws.Range("E34").Value = '  +1.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.102'
$ws.Range("E35").Value = '  +2.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.712'
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01835'
$ws.Range("E37").Value = '  +0.97%  '
$ws.Range("E38").Value = '  +1.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9226'
$ws.Range("E39").Value = '  +0.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.061'
$ws.Range("E40").Value = '  +1.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4417'
$ws.Range("E41").Value = '  +2.69%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.947'
$ws.Range("E42").Value = '  +4.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '106.41'
$ws.Range("E43").Value = '  +2.26%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9940'
$ws.Range("E44").Value = '  -0.56%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.33'
$ws.Range("E45").Value = '  +21.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1385'
$ws.Range("E46").Value = '  +3.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.611'
$ws.Range("E47").Value = '  +4.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.000'
$ws.Range("E48").Value = '  +4.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.97'
$ws.Range("E49").Value = '  +5.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05830'
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.426'
$ws.Range("E51").Value = '  +1.22%  '
